# Update evaluation metrics for the "RF" (row 6) and "Ensemble" (row 7) rows
# with newly (re-)computed precision/recall/f1/accuracy values per class/model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row6 = @(
    0.8305879661404714,
    0.8396972991405803,
    0.8305879661404714,
    0.8319988464678246,
    0.8112331274307939,
    0.8233575453282475,
    0.8112331274307939,
    0.8129002063521378,
    0.7253717684740335,
    0.7324787718191119,
    0.7253717684740335,
    0.7248313897793925,
    0.7962251201098148,
    0.807798954614249,
    0.7962251201098148,
    0.7967742207061101,
    0.7940288263555251,
    0.8017485226484183,
    0.7940288263555251,
    0.7949131375952623,
    0.8219629375428964,
    0.825416323160782,
    0.8219629375428964,
    0.8213291381248924
)

$row7 = @(
    0.8433539235872798,
    0.8487159913715232,
    0.8433539235872798,
    0.8438845433216151,
    0.8584305650880806,
    0.8606353889948695,
    0.8584305650880806,
    0.8573795959786967,
    0.8283687943262411,
    0.8396140992435329,
    0.8283687943262411,
    0.8293527408589704,
    0.8455044612216884,
    0.8497538216940527,
    0.8455044612216884,
    0.8455985765736322,
    0.8562113932738503,
    0.8647471614363464,
    0.8562113932738503,
    0.856295760917458,
    0.8369251887439946,
    0.8420709957133198,
    0.8369251887439946,
    0.837194762515152
)

for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, $i + 2).Value = $row6[$i]
}

for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, $i + 2).Value = $row7[$i]
}
